$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.1746478873239437
$ws.Range("C2").Value = 0.5859154929577465
$ws.Range("J2").Value = 0.005633802816901409
$ws.Range("P2").Value = 0.123943661971831
$ws.Range("S2").Value = 0.1098591549295775

# Row 3
$ws.Range("B3").Value = 0.004672897196261682
$ws.Range("C3").Value = 0.02336448598130841
$ws.Range("J3").Value = 0.02803738317757009
$ws.Range("P3").Value = 0.7616822429906542
$ws.Range("S3").Value = 0.1822429906542056

# Row 4
$ws.Range("J4").Value = 0.0392156862745098
$ws.Range("O4").Value = 0.0196078431372549
$ws.Range("P4").Value = 0.5490196078431373
$ws.Range("S4").Value = 0.392156862745098

# Row 6
$ws.Range("B6").Value = 0.09012875536480687
$ws.Range("D6").Value = 0.008583690987124463
$ws.Range("F6").Value = 0.0944206008583691
$ws.Range("J6").Value = 0.2746781115879828
$ws.Range("O6").Value = 0.02145922746781116
$ws.Range("Q6").Value = 0.2017167381974249
$ws.Range("R6").Value = 0.07296137339055794
$ws.Range("S6").Value = 0.2360515021459227

# Row 7
$ws.Range("B7").Value = 0.1038251366120219
$ws.Range("D7").Value = 0.02185792349726776
$ws.Range("F7").Value = 0.02185792349726776
$ws.Range("J7").Value = 0.1530054644808743
$ws.Range("O7").Value = 0.01092896174863388
$ws.Range("Q7").Value = 0.1639344262295082
$ws.Range("R7").Value = 0.1092896174863388
$ws.Range("S7").Value = 0.4153005464480874

# Row 8
$ws.Range("B8").Value = 0.1199095022624434
$ws.Range("D8").Value = 0.00904977375565611
$ws.Range("F8").Value = 0.05882352941176471
$ws.Range("J8").Value = 0.1357466063348416
$ws.Range("O8").Value = 0.01809954751131222
$ws.Range("Q8").Value = 0.1719457013574661
$ws.Range("R8").Value = 0.08823529411764706
$ws.Range("S8").Value = 0.3981900452488688

# Row 9
$ws.Range("B9").Value = 0.1065573770491803
$ws.Range("D9").Value = 0.01229508196721311
$ws.Range("F9").Value = 0.02868852459016394
$ws.Range("J9").Value = 0.1229508196721311
$ws.Range("O9").Value = 0.02868852459016394
$ws.Range("Q9").Value = 0.209016393442623
$ws.Range("R9").Value = 0.110655737704918
$ws.Range("S9").Value = 0.3811475409836065

# Row 10
$ws.Range("B10").Value = 0.1204225352112676
$ws.Range("D10").Value = 0.02887323943661972
$ws.Range("E10").Value = 0.002112676056338028
$ws.Range("F10").Value = 0.06267605633802817
$ws.Range("J10").Value = 0.1323943661971831
$ws.Range("O10").Value = 0.01267605633802817
$ws.Range("Q10").Value = 0.2133802816901408
$ws.Range("R10").Value = 0.0795774647887324
$ws.Range("S10").Value = 0.347887323943662

# Row 11
$ws.Range("G11").Value = 0.125
$ws.Range("J11").Value = 0.1118421052631579
$ws.Range("K11").Value = 0.1875
$ws.Range("L11").Value = 0.555921052631579
$ws.Range("S11").Value = 0.01973684210526316

# Row 12
$ws.Range("G12").Value = 0.7352941176470589
$ws.Range("J12").Value = 0.2058823529411765
$ws.Range("K12").Value = 0.02352941176470588
$ws.Range("L12").Value = 0.01764705882352941
$ws.Range("S12").Value = 0.01764705882352941

# Row 15
$ws.Range("F15").Value = 0.05
$ws.Range("H15").Value = 0.1291666666666667
$ws.Range("I15").Value = 0.0625
$ws.Range("J15").Value = 0.3625
$ws.Range("K15").Value = 0.0625
$ws.Range("O15").Value = 0.09583333333333334
$ws.Range("S15").Value = 0.2375

# Row 16
$ws.Range("F16").Value = 0.02597402597402598
$ws.Range("H16").Value = 0.1471861471861472
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.354978354978355
$ws.Range("K16").Value = 0.09956709956709957
$ws.Range("M16").Value = 0.01731601731601732
$ws.Range("O16").Value = 0.08225108225108226
$ws.Range("S16").Value = 0.1818181818181818

# Row 17
$ws.Range("F17").Value = 0.02376237623762376
$ws.Range("H17").Value = 0.1663366336633663
$ws.Range("I17").Value = 0.07920792079207921
$ws.Range("J17").Value = 0.4396039603960396
$ws.Range("K17").Value = 0.05742574257425743
$ws.Range("M17").Value = 0.02772277227722772
$ws.Range("N17").Value = 0.00396039603960396
$ws.Range("O17").Value = 0.08514851485148515
$ws.Range("S17").Value = 0.1168316831683168

# Row 18
$ws.Range("F18").Value = 0.0319634703196347
$ws.Range("H18").Value = 0.136986301369863
$ws.Range("I18").Value = 0.1552511415525114
$ws.Range("J18").Value = 0.3378995433789954
$ws.Range("K18").Value = 0.0958904109589041
$ws.Range("M18").Value = 0.0091324200913242
$ws.Range("O18").Value = 0.0730593607305936
$ws.Range("S18").Value = 0.1598173515981735

# Row 19
$ws.Range("F19").Value = 0.01408450704225352
$ws.Range("H19").Value = 0.1934766493699036
$ws.Range("I19").Value = 0.09933283914010378
$ws.Range("J19").Value = 0.3787991104521868
$ws.Range("K19").Value = 0.1089696071163825
$ws.Range("M19").Value = 0.01482579688658265
$ws.Range("O19").Value = 0.05707931801334322
$ws.Range("S19").Value = 0.1334321719792439
